$d = $word.ActiveDocument

$cr = [char]13
$replacement = "INTRODUCCIÓN" + $cr + $cr + "La necesidad de disponer de una gestión de los medicamentos y un sistema que ayude a las ventas conduce al uso de herramientas informáticas de gestión"

$find = $d.Content.Find
$ok = $find.Execute("INTRODUCCIÓN", $false, $false, $false, $false, $false, $true, 1, $false, $replacement, 2)
